$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Extension_Category")

# --- Row 2: Transport Model ---
$ws.Range("C2").Value = "Track transport of pollutants, particles, or objects"

# --- Row 3: Water Quality Model -> Acidification Model ---
$ws.Range("B3").Value = "Acidification Model"
$ws.Range("D3").Value = 2

# --- Row 6: Other ---
$ws.Range("D6").Value = 14

# --- Row 7: Flood Model ---
$ws.Range("C7").Value = "Predicts flooding or inundation"

# --- Row 8: Stormwater/Drainage Model -> Watershed Model ---
$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "Watershed Model"
$ws.Range("C8").Value = "Watershed model or urban drainage model of water flow and nutrient loads"
$ws.Range("D8").Value = 8

# --- Row 9: River Discharge Model -> Ecosystem Model ---
$ws.Range("A9").Value = 12
$ws.Range("B9").Value = "Ecosystem Model"
$ws.Range("C9").Value = "Predict ecosystem response to conditions and nutrient loading"
$ws.Range("D9").Value = 1

# --- Row 10: Watershed Model -> Ice models ---
$ws.Range("A10").Value = 13
$ws.Range("B10").Value = "Ice models"
$ws.Range("C10").Value = "Models that predict freezing of the harbor, ice damage, or icing of ship and boat superstructure"
$ws.Range("D10").Value = 10

# --- Row 11: Not Specified -> Habitat Models ---
$ws.Range("A11").Value = 16
$ws.Range("B11").Value = "Habitat Models"
$ws.Range("C11").Value = "Models that predict habitat conditions or suitability"

# --- Row 12: Ecosystem Model -> Use Capabilities ---
$ws.Range("A12").Value = 17
$ws.Range("B12").Value = "Use Capabilities"
$ws.Range("C12").Value = "Extensions that allow application of the model in different ways or settings"
$ws.Range("D12").Value = 12

# --- Row 13: New row, Model Refinement ---
$ws.Range("A13").Value = 18
$ws.Range("B13").Value = "Model Refinement"
$ws.Range("C13").Value = "Improvements in resolution or accuracy of the model to address specific needs"
$ws.Range("D13").Value = 13

# --- Update defined name to extend through row 13 ---
$wb.Names.Item("Extension_Category").RefersTo = '=''Extension_Category''!$A$1:$D$13'
